# adding scores for 17th may
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marksheet")

# Row 29 - Physics
$ws.Range("C29").Value = "Physics"
$ws.Range("D29").Value = 48
$ws.Range("E29").Value = 43
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 3

# Row 30 - Biology
$ws.Range("C30").Value = "Biology"
$ws.Range("D30").Value = 25
$ws.Range("E30").Value = 25
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0

# Row 31 - Chemistry
$ws.Range("C31").Value = "Chemistry"
$ws.Range("D31").Value = 50
$ws.Range("E31").Value = 48
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 1

$ws.Activate()
$ws.Range("G31").Select()
